$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 0.1.0 -> 0.2.0
$meta.Range("B3").Value = "0.2.0"

# Date: 2022-01-25T10:18:17-05:00 -> 2022-02-08T17:12:45-05:00
$meta.Range("B8").Value = "2022-02-08T17:12:45-05:00"

# Case Sensitive: (blank) -> "true"
# Count: 35 -> 36
# A leading apostrophe forces Excel to store these as literal text instead
# of auto-converting "true"/"36" into a boolean / number.
$meta.Range("B14").Value = "'true"
$meta.Range("B21").Value = "'36"

# The apostrophe trick marks the cells with a "quote prefix" style variant;
# re-apply the plain data-row formatting (as used by the rest of column B)
# from neighboring cells so B14/B21 end up with the same style as the rest
# of the table.
$meta.Range("B13").Copy()
$meta.Range("B14").PasteSpecial(-4122) | Out-Null
$meta.Range("B20").Copy()
$meta.Range("B21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Concepts sheet: append the new "Lab Interpretive Report" concept row
# ---------------------------------------------------------------------
$concepts = $wb.Worksheets.Item("Concepts")

# Duplicate the formatting of the last existing data row (36) onto the new
# row 37 first (this also creates the empty, styled D37 cell).
$concepts.Range("A36:D36").Copy()
$concepts.Range("A37:D37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new row's values. Column A ("Level") is text ("1") throughout
# the table, so force it with a leading apostrophe rather than letting it
# become a number.
$concepts.Cells.Item(37, 1).Value = "'1"
$concepts.Cells.Item(37, 2).Value = "lab-interpretative-report"
$concepts.Cells.Item(37, 3).Value = "Lab Interpretive Report"

# Re-normalize A37's formatting (the apostrophe again leaves a quote-prefix
# style) back to the shared data-row style used across the column.
$concepts.Range("C36").Copy()
$concepts.Range("A37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
